$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold plain text (prices/percentages formatted as strings, e.g.
# "64.327.91" or "  -0.04%  "). Force Text number format on the affected range first
# so Excel does not reinterpret these strings as numbers (which would introduce
# floating point rounding and change the cell type), then restore the default
# "Normal" style so no stray formatting is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '64.327.91'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '3.407.58'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '569.79'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").Value = '156.38'
$ws.Range("E6").Value = '  -2.06%  '
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  +8.21%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '3.407.20'
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("D10").Value = '7.14'
$ws.Range("E10").Value = '  -2.16%  '
$ws.Range("E11").Value = '  -1.32%  '
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("D13").Value = '3.992.52'
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("E15").Value = '  -2.15%  '
$ws.Range("D16").Value = '27.51'
$ws.Range("E16").Value = '  -1.50%  '
$ws.Range("D17").Value = '64.290.51'
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").Value = '3.407.14'
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("D20").Value = '13.78'
$ws.Range("E20").Value = '  -2.43%  '
$ws.Range("D21").Value = '377.34'
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").Value = '8.01'
$ws.Range("E22").Value = '  -1.08%  '
$ws.Range("E23").Value = '  +0.73%  '
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = '71.50'
$ws.Range("E25").Value = '  -1.70%  '
$ws.Range("E26").Value = '  -1.99%  '
$ws.Range("E27").Value = '  +9.06%  '
$ws.Range("D28").Value = '0.178'
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.36%  '
$ws.Range("E30").Value = '  +4.62%  '
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").Value = '23.04'
$ws.Range("E33").Value = '  -1.73%  '
$ws.Range("D34").Value = '7.15'
$ws.Range("E34").Value = '  +1.40%  '
$ws.Range("E35").Value = '  +7.92%  '
$ws.Range("D36").Value = '159.67'
$ws.Range("E36").Value = '  -1.52%  '
$ws.Range("D37").Value = '1.89'
$ws.Range("E37").Value = '  +0.84%  '
$ws.Range("D38").Value = '6.95'
$ws.Range("E38").Value = '  +6.45%  '
$ws.Range("D39").Value = '0.0760'
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '2.877.30'
$ws.Range("E40").Value = '  -4.89%  '
$ws.Range("D41").Value = '4.62'
$ws.Range("E41").Value = '  +2.70%  '
$ws.Range("D42").Value = '26.26'
$ws.Range("E42").Value = '  -3.29%  '
$ws.Range("D43").Value = '42.94'
$ws.Range("E43").Value = '  +1.02%  '
$ws.Range("D44").Value = '0.0314'
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("D45").Value = '26.07'
$ws.Range("E45").Value = '  +6.92%  '
$ws.Range("D46").Value = '0.766'
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").Value = '321.79'
$ws.Range("E47").Value = '  +8.01%  '
$ws.Range("D48").Value = '1.07'
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("E49").Value = '  +2.98%  '
$ws.Range("E50").Value = '  +1.07%  '
$ws.Range("D51").Value = '6.55'
$ws.Range("E51").Value = '  -0.71%  '

$dataRange.Style = "Normal"
